# Update the "想去人数" (interested-count) figures in the 展览 and 全部类型
# sheets to reflect the latest generated output (gh-pages commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1351
$ws1.Range("F7").Value  = 11659
$ws1.Range("F8").Value  = 4381
$ws1.Range("F14").Value = 1093
$ws1.Range("F17").Value = 5080
$ws1.Range("F19").Value = 183
$ws1.Range("F20").Value = 513
$ws1.Range("F22").Value = 11261

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1351
$ws4.Range("F7").Value  = 11659
$ws4.Range("F8").Value  = 4381
$ws4.Range("F15").Value = 1093
$ws4.Range("F18").Value = 5080
$ws4.Range("F20").Value = 183
$ws4.Range("F21").Value = 513
$ws4.Range("F23").Value = 11261
